# Update cryptocurrency price (D) and 1h volume-change (E) figures
# per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.052.06"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "'3.201.90"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'574.29"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'167.08"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").Value = "'0.596"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "'0.391"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'3.763.45"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'65.088.59"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "'25.55"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'3.210.59"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "'0.0000157"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "'410.94"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'69.67"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("D25").Value = "'0.491"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -6.15%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "'1.84"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'21.54"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").Value = "'4.96"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'6.40"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'156.82"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("D36").Value = "'2.753.25"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'1.72"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'24.21"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("D39").Value = "'4.14"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'0.716"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "'0.0634"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'5.66"
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'296.56"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("D45").Value = "'21.45"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "'1.95"
$ws.Range("E48").Value = "  -10.46%  "
$ws.Range("D49").Value = "'5.79"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'0.902"
$ws.Range("E51").Value = "  -2.85%  "
